$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B64").NumberFormat = "@"

$ws.Range("B2").Value = "1010010000000000100000100000000010000110000000001000000000000000"
$ws.Range("D2").Value = 10
$ws.Range("B3").Value = "1110010000000000100000100000000010000110000000001000000000000000"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 9
$ws.Range("B4").Value = "1010000000000000100000000000000000000000000000000000000000000000"
$ws.Range("C4").Value = 0.015632
$ws.Range("D4").Value = 10
$ws.Range("B5").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C5").Value = 0
$ws.Range("B6").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C6").Value = 0.015629
$ws.Range("B7").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B8").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B9").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B10").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B11").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B12").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B13").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B14").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B15").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B16").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B17").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C17").Value = 0
$ws.Range("B18").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C18").Value = 0
$ws.Range("B19").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B20").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B21").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B22").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C22").Value = 0.015608
$ws.Range("B23").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B24").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B25").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B26").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B27").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B28").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B29").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B30").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C30").Value = 0
$ws.Range("B31").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C31").Value = 0
$ws.Range("B32").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C32").Value = 0
$ws.Range("B33").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C33").Value = 0
$ws.Range("B34").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C34").Value = 0
$ws.Range("B35").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C35").Value = 0.015624
$ws.Range("B36").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C36").Value = 0
$ws.Range("B37").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C37").Value = 0
$ws.Range("B38").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C38").Value = 0
$ws.Range("B39").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C39").Value = 0
$ws.Range("B40").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C40").Value = 0
$ws.Range("B41").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C41").Value = 0
$ws.Range("B42").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B43").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B44").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C44").Value = 0.012052
$ws.Range("B45").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C45").Value = 0.00187
$ws.Range("B46").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B47").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B48").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C48").Value = 0.0072
$ws.Range("B49").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C49").Value = 0.004283
$ws.Range("B50").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C50").Value = 0.004997
$ws.Range("B51").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C51").Value = 0.001
$ws.Range("B52").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C52").Value = 0
$ws.Range("B53").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B54").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C54").Value = 0.015714
$ws.Range("B55").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("B56").Value = "1110101100000000101000001000000010000100000000001000001000000000"
$ws.Range("C56").Value = 0
$ws.Range("B57").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("C57").Value = 0.015626
$ws.Range("D57").Value = 1
$ws.Range("B58").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("D58").Value = 1
$ws.Range("B59").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 1
$ws.Range("B60").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("C60").Value = 0.015624
$ws.Range("D60").Value = 1
$ws.Range("B61").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("D61").Value = 1
$ws.Range("B62").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 1
$ws.Range("B63").Value = "1110101100000000101000001000000010000000000000001000001000000000"
$ws.Range("D63").Value = 1
$ws.Range("B64").Value = "1110101100000000101000001000000010000000000000001000000000000000"
$ws.Range("C64").Value = 0.015627
$ws.Range("D64").Value = 2
